# Update on 7/21/2025 at 2:29pm
# - Rework the FY26 (Estimate) row (row 8) of the "Table" sheet so that:
#     F8 (Total Enrollee Months) is now derived from the YTD figures (F9/(B9/12))
#     G8 (Avg Annual Enrollment) is now an average of F8 and F10 (monthly)
#     H8 (End of Period Enrollment) switches from "NA" to a computed value (F8/12)
#     J8 (Enrollee Months % Growth) is rewritten as (F8-F10)/F10
#   Downstream formulas (rows 4-7) are plain formula cells that recalculate
#   automatically once these inputs change.
# - J7 is rewritten to simply reference the growth-rate input cell H25.
# - H9:H20 (End of Period Enrollment, actuals) switch from the old "NA/accounting"
#   number format to the plain #,##0 number format used elsewhere in the table.
# - The "FY27 Growth Assumptions" label (H22) is renamed to "FY27-FY30 Growth
#   Assumptions" to reflect the extended forecast horizon.
# - The active selection ends up on G8, the last cell touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table")

# --- Row 8 (FY26 (Estimate)) formula rewrites -----------------------------

# F8: Total Enrollee Months now computed from YTD months reported (F9 / (B9/12))
$ws.Range("F8").Formula = "=F9/(B9/12)"

# G8: Average Annual Enrollment now the average of F8 and F10, monthly
$ws.Range("G8").Formula = "=((F8+F10)/2)/12"

# H8: End of Period Enrollment - was a literal "NA" text cell (accounting-style
# number format); becomes a real formula with the plain #,##0 number format
# that the rest of column H/"Average Annual Enrollment" columns use. Copy the
# number format from a cell that already carries that exact style (C4) before
# setting the formula so the underlying style entry is reused rather than a
# near-duplicate style being minted.
$ws.Range("C4").Copy()
$ws.Range("H8").PasteSpecial(-4122)
$ws.Range("H8").Formula = "=F8/12"

# J8: Enrollee Months % Growth - rewritten to a direct percentage-difference
# formula and restyled with the 0.0% percent format (style used by J10:J19).
$ws.Range("J10").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J8").Formula = "=(F8-F10)/F10"

# --- Row 7 tweak ------------------------------------------------------------

# J7: now simply mirrors the growth-rate input cell H25 instead of being
# computed from I7/F8.
$ws.Range("J7").Formula = "=H25"

# --- Row 6 formula (kept identical text, no longer part of the shared group) -

$ws.Range("J6").Formula = "=I6/F7"

# --- H9:H20 number-format change (accounting -> plain #,##0) ---------------

$ws.Range("C4").Copy()
$ws.Range("H9:H20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Label rename: "FY27 Growth Assumptions" -> "FY27-FY30 Growth Assumptions"

$ws.Range("H22").Value = "FY27-FY30 Growth Assumptions"

# --- Final selection on G8, matching the last-edited cell ------------------

$ws.Activate()
$ws.Range("G8").Select()
